$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2262295081967213
$ws.Range("C2").Value = 0.5016393442622951
$ws.Range("P2").Value = 0.180327868852459
$ws.Range("S2").Value = 0.09180327868852459
$ws.Range("B3").Value = 0.0130718954248366
$ws.Range("C3").Value = 0.03267973856209151
$ws.Range("J3").Value = 0.006535947712418301
$ws.Range("O3").Value = 0.006535947712418301
$ws.Range("P3").Value = 0.7320261437908496
$ws.Range("S3").Value = 0.2091503267973856
$ws.Range("J4").Value = 0.04081632653061224
$ws.Range("P4").Value = 0.6326530612244898
$ws.Range("S4").Value = 0.3265306122448979
$ws.Range("B6").Value = 0.06403940886699508
$ws.Range("D6").Value = 0.03448275862068965
$ws.Range("F6").Value = 0.03940886699507389
$ws.Range("J6").Value = 0.2512315270935961
$ws.Range("O6").Value = 0.01477832512315271
$ws.Range("Q6").Value = 0.1724137931034483
$ws.Range("R6").Value = 0.05911330049261083
$ws.Range("S6").Value = 0.3645320197044335
$ws.Range("B7").Value = 0.08602150537634409
$ws.Range("D7").Value = 0.04301075268817205
$ws.Range("F7").Value = 0.05376344086021505
$ws.Range("J7").Value = 0.1827956989247312
$ws.Range("O7").Value = 0.005376344086021506
$ws.Range("Q7").Value = 0.1182795698924731
$ws.Range("R7").Value = 0.05913978494623656
$ws.Range("S7").Value = 0.4516129032258064
$ws.Range("B8").Value = 0.09392265193370165
$ws.Range("D8").Value = 0.01657458563535912
$ws.Range("F8").Value = 0.06998158379373849
$ws.Range("J8").Value = 0.1270718232044199
$ws.Range("O8").Value = 0.01657458563535912
$ws.Range("Q8").Value = 0.1583793738489871
$ws.Range("R8").Value = 0.09576427255985268
$ws.Range("S8").Value = 0.4217311233885819
$ws.Range("B9").Value = 0.07462686567164178
$ws.Range("D9").Value = 0.009950248756218905
$ws.Range("E9").Value = 0.004975124378109453
$ws.Range("F9").Value = 0.05472636815920398
$ws.Range("J9").Value = 0.1592039800995025
$ws.Range("O9").Value = 0.009950248756218905
$ws.Range("Q9").Value = 0.1492537313432836
$ws.Range("R9").Value = 0.08955223880597014
$ws.Range("S9").Value = 0.4477611940298508
$ws.Range("B10").Value = 0.1138732959101844
$ws.Range("D10").Value = 0.02085004009623095
$ws.Range("F10").Value = 0.06174819566960706
$ws.Range("J10").Value = 0.1299117882919006
$ws.Range("O10").Value = 0.009623095429029671
$ws.Range("Q10").Value = 0.1924619085805934
$ws.Range("R10").Value = 0.07457898957497995
$ws.Range("S10").Value = 0.3969526864474739
$ws.Range("G11").Value = 0.1590909090909091
$ws.Range("K11").Value = 0.2402597402597403
$ws.Range("L11").Value = 0.474025974025974
$ws.Range("S11").Value = 0.01948051948051948
$ws.Range("G12").Value = 0.7597402597402597
$ws.Range("J12").Value = 0.1428571428571428
$ws.Range("K12").Value = 0.006493506493506494
$ws.Range("L12").Value = 0.03896103896103896
$ws.Range("S12").Value = 0.05194805194805195
$ws.Range("G13").Value = 0.6097560975609756
$ws.Range("J13").Value = 0.3414634146341464
$ws.Range("S13").Value = 0.04878048780487805
$ws.Range("F15").Value = 0.01639344262295082
$ws.Range("H15").Value = 0.1912568306010929
$ws.Range("I15").Value = 0.08743169398907104
$ws.Range("J15").Value = 0.3387978142076503
$ws.Range("K15").Value = 0.06010928961748634
$ws.Range("M15").Value = 0.00546448087431694
$ws.Range("O15").Value = 0.01639344262295082
$ws.Range("S15").Value = 0.2841530054644809
$ws.Range("F16").Value = 0.02072538860103627
$ws.Range("H16").Value = 0.1813471502590674
$ws.Range("I16").Value = 0.1243523316062176
$ws.Range("J16").Value = 0.3471502590673575
$ws.Range("K16").Value = 0.1088082901554404
$ws.Range("O16").Value = 0.04145077720207254
$ws.Range("S16").Value = 0.1761658031088083
$ws.Range("F17").Value = 0.01941747572815534
$ws.Range("H17").Value = 0.1868932038834951
$ws.Range("I17").Value = 0.09223300970873786
$ws.Range("J17").Value = 0.3810679611650485
$ws.Range("K17").Value = 0.0970873786407767
$ws.Range("M17").Value = 0.01456310679611651
$ws.Range("N17").Value = 0.002427184466019417
$ws.Range("O17").Value = 0.04854368932038835
$ws.Range("S17").Value = 0.1577669902912621
$ws.Range("F18").Value = 0.01081081081081081
$ws.Range("H18").Value = 0.1891891891891892
$ws.Range("I18").Value = 0.0918918918918919
$ws.Range("J18").Value = 0.4054054054054054
$ws.Range("K18").Value = 0.05945945945945946
$ws.Range("M18").Value = 0.01621621621621622
$ws.Range("O18").Value = 0.01081081081081081
$ws.Range("S18").Value = 0.2162162162162162
$ws.Range("F19").Value = 0.01227830832196453
$ws.Range("H19").Value = 0.2482946793997272
$ws.Range("I19").Value = 0.07503410641200546
$ws.Range("J19").Value = 0.3287858117326057
$ws.Range("K19").Value = 0.09959072305593451
$ws.Range("M19").Value = 0.02046384720327422
$ws.Range("N19").Value = 0.0006821282401091405
$ws.Range("O19").Value = 0.06343792633015007
$ws.Range("S19").Value = 0.1514324693042292
